# Add "2022-Q4" sheet (new quarter of fund-holdings data) right after the
# "总计" (totals) sheet, in front of "2022-Q3", and record the new quarter
# in the totals sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Create the new "2022-Q4" worksheet.
#    Easiest reliable way to get identical styling (header row + index
#    column formatting) is to duplicate the existing "2022-Q3" sheet
#    (same shape: header + 5 fund rows) and then overwrite its values.
# ------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2022-Q3")
$srcSheet.Copy($srcSheet) | Out-Null
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# Row 2 : 009649 - 嘉实精选平衡混合A
$q4.Range("B2").Value = "'009649"
$q4.Range("C2").Value = "嘉实精选平衡混合A"
$q4.Range("D2").Value = "'0.47"
$q4.Range("E2").Value = "'68.05"
$q4.Range("F2").Value = "'2.72"
$q4.Range("G2").Value = "'0.0128"
$q4.Range("H2").Value = 10

# Row 3 : 001563 - 华富健康文娱灵活配置混合
$q4.Range("B3").Value = "'001563"
$q4.Range("C3").Value = "华富健康文娱灵活配置混合"
$q4.Range("D3").Value = "'0.33"
$q4.Range("E3").Value = "'90.72"
$q4.Range("F3").Value = "'3.65"
$q4.Range("G3").Value = "'0.0120"
$q4.Range("H3").Value = 8

# Row 4 : 008884 - 博远博锐混合A
$q4.Range("B4").Value = "'008884"
$q4.Range("C4").Value = "博远博锐混合A"
$q4.Range("D4").Value = "'0.12"
$q4.Range("E4").Value = "'73.80"
$q4.Range("F4").Value = "'2.53"
$q4.Range("G4").Value = "'0.0030"
$q4.Range("H4").Value = 8

# Row 5 : 009650 - 嘉实精选平衡混合C
$q4.Range("B5").Value = "'009650"
$q4.Range("C5").Value = "嘉实精选平衡混合C"
$q4.Range("D5").Value = "'0.04"
$q4.Range("E5").Value = "'68.05"
$q4.Range("F5").Value = "'2.72"
$q4.Range("G5").Value = "'0.0011"
$q4.Range("H5").Value = 10

# Row 6 : 008885 - 博远博锐混合C
$q4.Range("B6").Value = "'008885"
$q4.Range("C6").Value = "博远博锐混合C"
$q4.Range("D6").Value = "'0.02"
$q4.Range("E6").Value = "'73.80"
$q4.Range("F6").Value = "'2.53"
$q4.Range("G6").Value = "'0.0005"
$q4.Range("H6").Value = 8

# ------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert a new row right after the
#    header for 2022-Q4, shifting the older quarters down, and bump the
#    running index in column A.
# ------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Rows("2:2").Insert() | Out-Null
$totals.Range("B2:D2").ClearFormats() | Out-Null

$totals.Range("A3").Copy() | Out-Null
$totals.Range("A2").PasteSpecial(-4122) | Out-Null

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q4"
$totals.Range("C2").Value = 5
$totals.Range("D2").Value = 0.03

$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4
$totals.Range("A7").Value = 5
